$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.292.60'
$ws.Range("E2").Value = '  +0.13%  '
$ws.Range("D3").Value = '1.689.49'
$ws.Range("E3").Value = '  +1.04%  '
$ws.Range("E4").Value = '  +0.29%  '
$ws.Range("E5").Value = '  +0.50%  '
$ws.Range("D6").Value = '''0.5252'
$ws.Range("E6").Value = '  +2.63%  '
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("D8").Value = '''0.2700'
$ws.Range("E8").Value = '  +1.37%  '
$ws.Range("D9").Value = '''0.06438'
$ws.Range("E9").Value = '  +1.15%  '
$ws.Range("E10").Value = '  +1.60%  '
$ws.Range("D11").Value = '''0.07461'
$ws.Range("E11").Value = '  +1.30%  '
$ws.Range("D12").Value = '1.697.47'
$ws.Range("E12").Value = '  +1.53%  '
$ws.Range("E13").Value = '  +0.03%  '
$ws.Range("D14").Value = '''0.5860'
$ws.Range("E14").Value = '  +0.83%  '
$ws.Range("D15").Value = '''0.000008527'
$ws.Range("E15").Value = '  -0.62%  '
$ws.Range("D16").Value = '''64.62'
$ws.Range("E16").Value = '  -0.86%  '
$ws.Range("D17").Value = '26.332.56'
$ws.Range("E17").Value = '  +0.10%  '
$ws.Range("D18").Value = '''4.959'
$ws.Range("E18").Value = '  +0.20%  '
$ws.Range("D19").Value = '''1.007'
$ws.Range("E19").Value = '  +0.19%  '
$ws.Range("E20").Value = '  +0.20%  '
$ws.Range("D21").Value = '''189.73'
$ws.Range("E21").Value = '  +0.00%  '
$ws.Range("D22").Value = '''6.222'
$ws.Range("E22").Value = '  +0.23%  '
$ws.Range("D23").Value = '''1.008'
$ws.Range("E23").Value = '  +0.29%  '
$ws.Range("D24").Value = '''144.80'
$ws.Range("E24").Value = '  +0.90%  '
$ws.Range("D25").Value = '''7.671'
$ws.Range("E25").Value = '  -0.14%  '
$ws.Range("E26").Value = '  +4.84%  '
$ws.Range("E27").Value = '  +0.54%  '
$ws.Range("D28").Value = '''0.06686'
$ws.Range("E28").Value = '  +15.08%  '
$ws.Range("D29").Value = '''1.349'
$ws.Range("E29").Value = '  +5.06%  '
$ws.Range("E30").Value = '  +0.79%  '
$ws.Range("D31").Value = '''3.588'
$ws.Range("E31").Value = '  +1.99%  '
$ws.Range("E32").Value = '  +0.99%  '
$ws.Range("D33").Value = '''1.671'
$ws.Range("E33").Value = '  +1.05%  '
$ws.Range("E34").Value = '  +1.94%  '
$ws.Range("D35").Value = '''0.6222'
$ws.Range("E35").Value = '  +3.72%  '
$ws.Range("D36").Value = '''2.393'
$ws.Range("E36").Value = '  +1.61%  '
$ws.Range("D37").Value = '''2.691'
$ws.Range("E37").Value = '  +1.85%  '
$ws.Range("D38").Value = '''6.337'
$ws.Range("E38").Value = '  +5.49%  '
$ws.Range("D39").Value = '''0.01624'
$ws.Range("E39").Value = '  +0.55%  '
$ws.Range("D40").Value = '1.106.26'
$ws.Range("E40").Value = '  +2.23%  '
$ws.Range("D41").Value = '''0.8852'
$ws.Range("E41").Value = '  +3.05%  '
$ws.Range("E42").Value = '  +0.82%  '
$ws.Range("D43").Value = '''101.34'
$ws.Range("E43").Value = '  +1.47%  '
$ws.Range("D44").Value = '1.837.32'
$ws.Range("E44").Value = '  +1.09%  '
$ws.Range("E45").Value = '  +1.14%  '
$ws.Range("D46").Value = '''56.86'
$ws.Range("E46").Value = '  +1.62%  '
$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").Value = '''1.009'
$ws.Range("E47").Value = '  +0.04%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '''8.170'
$ws.Range("E48").Value = '  +0.85%  '
$ws.Range("D49").Value = '''0.05265'
$ws.Range("E49").Value = '  +1.57%  '
$ws.Range("E50").Value = '  +0.23%  '
$ws.Range("D51").Value = '''6.059'
$ws.Range("E51").Value = '  +3.05%  '
